# Regenerate merged AHB files
# - Rename the _old/_new header-suffix columns to _FV2210/_FV2304
# - Wrap the data range in an Excel Table (ListObject)
# - Freeze the header row (pane split at row 2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row text (A1:U1) in place: _old -> _FV2210, _new -> _FV2304
$headers = @(
    "Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210",
    "Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
    "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the used range A1:U79 into an Excel Table named Table1, with banded rows
#    and no explicit table style (closest match to the source formatting).
$rng = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3) Freeze panes at row 1 (so row 1 stays visible while scrolling)
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true

Write-Host "done"
